$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7..36 down to 8..37
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new data point
$ws.Cells.Item(7, 1).Value = 1
$ws.Cells.Item(7, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(7, 4).Value = 44972
$ws.Cells.Item(7, 4).NumberFormat = $ws.Cells.Item(8, 4).NumberFormat
$ws.Cells.Item(7, 5).Value = 15
$ws.Cells.Item(7, 6).Value = 100112044
$ws.Cells.Item(7, 7).Value = "Perejil"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 350
$ws.Cells.Item(7, 11).Value = 800
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = 943
$ws.Cells.Item(7, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 472
$ws.Cells.Item(7, 17).Value = 2
$ws.Cells.Item(7, 18).Value = "Hortaliza"
